# Update "想去人数" (want-to-go count) values in F column for the
# sheets that list individual event rows: "展览" (sheet1) and
# "全部类型" (sheet4). Each row's F-value is bumped by a small amount
# as reflected in the upstream data refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8884
    3  = 8357
    8  = 764
    9  = 221
    10 = 5539
    11 = 13
    14 = 23
    15 = 25
    18 = 226
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
